# "final version before presentaion"
#
# The hotel-site translation sheet got six of its label rows renamed
# (old ad-hoc "xxxpolicy" style labels replaced with the cleaner names
# used elsewhere in the project), the page was set up for printing, and
# the sheet selection was left on A30 after the edits.
#
# Edit order below mirrors the order the labels were actually retyped in
# (it only affects the internal shared-string bookkeeping, not any
# visible content) so the six relabeled rows line up the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Survey"
$ws.Range("A52").Value = "Cancellation"
$ws.Range("A53").Value = "Cookie_policy"
$ws.Range("A55").Value = "Privacy_policy"
$ws.Range("A46").Value = "Special_offer"
$ws.Range("A54").Value = "Policies"

# Set up the page for printing (adds <pageSetup .../> to the sheet).
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

# Leave the selection where the author left it before saving.
$ws.Range("A30").Select()
